{"js": "// Add files via upload\n//\n// The author's edit adds a clause to the bullet point that ends with\n// \"...a v\u00e9rszeg\u00e9nys\u00e9get korrig\u00e1lni kell\" \u2014 appending:\n//   \", az orrv\u00e9rz\u00e9sek el\u0151fordul\u00e1s\u00e1t limit\u00e1lni a t\u00falzott v\u00e9rnyom\u00e1semelked\u00e9s\n//    ker\u00fcl\u00e9s\u00e9vel\"\n// so the bullet reads in full:\n//   \"S\u00falyos orrv\u00e9rz\u00e9s (epistaxis) vagy kr\u00f3nikus anaemia  a v\u00e9rszeg\u00e9nys\u00e9get\n//    korrig\u00e1lni kell, az orrv\u00e9rz\u00e9sek el\u0151fordul\u00e1s\u00e1t limit\u00e1lni a t\u00falzott\n//    v\u00e9rnyom\u00e1semelked\u00e9s ker\u00fcl\u00e9s\u00e9vel\"\n\nconst body = context.document.body;\n\nconst searchResults = body.search(\"v\u00e9rszeg\u00e9nys\u00e9get korrig\u00e1lni kell\", {\n  matchCase: false,\n  matchWholeWord: false,\n});\nsearchResults.load(\"items\");\nawait context.sync();\n\nif (searchResults.items.length === 0) {\n  throw new Error(\"Anchor phrase 'v\u00e9rszeg\u00e9nys\u00e9get korrig\u00e1lni kell' not found\");\n}\n\n// Use the last match in case the phrase appears more than once; insert right\n// after it, inheriting its formatting (Times New Roman, 14pt).\nconst anchor = searchResults.items[searchResults.items.length - 1];\nanchor.insertText(\n  \", az orrv\u00e9rz\u00e9sek el\u0151fordul\u00e1s\u00e1t limit\u00e1lni a t\u00falzott v\u00e9rnyom\u00e1semelked\u00e9s ker\u00fcl\u00e9s\u00e9vel\",\n  Word.InsertLocation.end\n);\n\nawait context.sync();\n", "ps1": "# Add files via upload\n#\n# The author's edit adds a clause to the bullet point that ends with\n# \"...a v\u00e9rszeg\u00e9nys\u00e9get korrig\u00e1lni kell\" -- appending:\n#   \", az orrv\u00e9rz\u00e9sek el\u0151fordul\u00e1s\u00e1t limit\u00e1lni a t\u00falzott v\u00e9rnyom\u00e1semelked\u00e9s\n#    ker\u00fcl\u00e9s\u00e9vel\"\n# so the bullet reads in full:\n#   \"S\u00falyos orrv\u00e9rz\u00e9s (epistaxis) vagy kr\u00f3nikus anaemia  a v\u00e9rszeg\u00e9nys\u00e9get\n#    korrig\u00e1lni kell, az orrv\u00e9rz\u00e9sek el\u0151fordul\u00e1s\u00e1t limit\u00e1lni a t\u00falzott\n#    v\u00e9rnyom\u00e1semelked\u00e9s ker\u00fcl\u00e9s\u00e9vel\"\n\n$d = $word.ActiveDocument\n\n$anchorText = \"v\u00e9rszeg\u00e9nys\u00e9get korrig\u00e1lni kell\"\n$addition = \", az orrv\u00e9rz\u00e9sek el\u0151fordul\u00e1s\u00e1t limit\u00e1lni a t\u00falzott v\u00e9rnyom\u00e1semelked\u00e9s ker\u00fcl\u00e9s\u00e9vel\"\n\n$rng = $d.Content\n$found = $rng.Find.Execute($anchorText)\n\nif ($found) {\n    # Appending directly to the matched Range's .Text keeps the inserted\n    # text inside the same run, inheriting its run formatting (font,\n    # size, language, etc.) instead of creating an unformatted run.\n    $rng.Text = $rng.Text + $addition\n} else {\n    Write-Output \"Anchor phrase not found: $anchorText\"\n}\n"}
